$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 971
$ws.Range("C3").Value = 1925
$ws.Range("D3").Value = 3680
$ws.Range("E3").Value = 8596
$ws.Range("F3").Value = 11100
$ws.Range("G3").Value = 11900

$ws.Range("B8").Value = 8082
$ws.Range("C8").Value = 12900
$ws.Range("D8").Value = 27800
$ws.Range("E8").Value = 58500
$ws.Range("F8").Value = 92700
$ws.Range("G8").Value = 125000

$ws.Range("B13").Value = 7670
$ws.Range("C13").Value = 11900
$ws.Range("D13").Value = 16200
$ws.Range("E13").Value = 19100
$ws.Range("F13").Value = 18300
$ws.Range("G13").Value = 18100

$ws.Range("B18").Value = 110000
$ws.Range("C18").Value = 189000
$ws.Range("D18").Value = 286000
$ws.Range("E18").Value = 365000
$ws.Range("F18").Value = 406000
$ws.Range("G18").Value = 465000

$ws.Range("B23").Value = 2730
$ws.Range("C23").Value = 4506
$ws.Range("D23").Value = 7816
$ws.Range("E23").Value = 9189
$ws.Range("F23").Value = 8094
$ws.Range("G23").Value = 7440

$ws.Range("B28").Value = 123000
$ws.Range("C28").Value = 151000
$ws.Range("D28").Value = 255000
$ws.Range("E28").Value = 277000
$ws.Range("F28").Value = 367000
$ws.Range("G28").Value = 256000

$ws.Range("B33").Value = 6872
$ws.Range("C33").Value = 9846
$ws.Range("D33").Value = 11700
$ws.Range("E33").Value = 12900
$ws.Range("F33").Value = 13400
$ws.Range("G33").Value = 14200

$ws.Range("B38").Value = 106000
$ws.Range("C38").Value = 192000
$ws.Range("D38").Value = 237000
$ws.Range("E38").Value = 302000
$ws.Range("F38").Value = 359000
$ws.Range("G38").Value = 414000
